$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 1353
$ws.Range("K3").Value = 1287
$ws.Range("C4").Value = 1845
$ws.Range("J4").Value = 1792
$ws.Range("K4").Value = 282
$ws.Range("K5").Value = 80
$ws.Range("K6").Value = 1648
$ws.Range("C7").Value = 28389
$ws.Range("J7").Value = 29262
$ws.Range("K7").Value = 4650

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 30
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 137
$ws.Range("K8").Value = 269
$ws.Range("K14").Value = 30
$ws.Range("K17").Value = 6
$ws.Range("K23").Value = 44
$ws.Range("K29").Value = 217
$ws.Range("K31").Value = 51
$ws.Range("J33").Value = 1317
$ws.Range("K33").Value = 191
$ws.Range("K37").Value = 157
$ws.Range("K41").Value = 49
$ws.Range("J42").Value = 1237
$ws.Range("K42").Value = 157
$ws.Range("K49").Value = 31
$ws.Range("K53").Value = 64
$ws.Range("K54").Value = 80
$ws.Range("K55").Value = 49
$ws.Range("K60").Value = 33
$ws.Range("K61").Value = 8
$ws.Range("C63").Value = 274
$ws.Range("K65").Value = 119
$ws.Range("K67").Value = 187
$ws.Range("K69").Value = 16
$ws.Range("K73").Value = 45
$ws.Range("K75").Value = 17
$ws.Range("K76").Value = 64
$ws.Range("K79").Value = 125
$ws.Range("K80").Value = 17
$ws.Range("K83").Value = 95
$ws.Range("K85").Value = 239
$ws.Range("K86").Value = 33
$ws.Range("K88").Value = 58
$ws.Range("K92").Value = 19
$ws.Range("K94").Value = 57
$ws.Range("K97").Value = 40
$ws.Range("K98").Value = 29
$ws.Range("C101").Value = 28389
$ws.Range("J101").Value = 29262
$ws.Range("K101").Value = 4650

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K2").Value = 13
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 51
$ws.Range("K3").Value = 41
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 137

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 32
$ws.Range("K3").Value = 22

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 86
$ws.Range("K3").Value = 77
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 239

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 16

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 80
$ws.Range("K3").Value = 80
$ws.Range("K6").Value = 88
$ws.Range("K7").Value = 269

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 43
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 95

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J4").Value = 57
$ws.Range("K6").Value = 49
$ws.Range("J7").Value = 1317
$ws.Range("K7").Value = 191

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 31
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 157

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 31
$ws.Range("K3").Value = 31
$ws.Range("K6").Value = 54
$ws.Range("K7").Value = 119

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 55
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 187

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K4").Value = 4
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 56
$ws.Range("K6").Value = 77
$ws.Range("K7").Value = 217

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K3").Value = 46
$ws.Range("J4").Value = 55
$ws.Range("K4").Value = 10
$ws.Range("J7").Value = 1237
$ws.Range("K7").Value = 157

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 46
$ws.Range("K3").Value = 42
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("K2").Value = 4
$ws.Range("K7").Value = 6

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K2").Value = 18
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 29

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K3").Value = 11
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("K2").Value = 4
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 19

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("K4").Value = 4
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K3").Value = 14
$ws.Range("K4").Value = 5
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("K6").Value = 8
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 8
